$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New column headers (row 2) ----
$ws.Range("D2").Value = "Название"
$ws.Range("E2").Value = "Управляющий сигнал с МП"

# ---- New column D values (pin "Название") ----
$ws.Range("D3").Value  = "А1"
$ws.Range("D4").Value  = "А2"
$ws.Range("D5").Value  = "А3"
$ws.Range("D6").Value  = "А4"
$ws.Range("D7").Value  = "А5"
$ws.Range("D8").Value  = "А6"
$ws.Range("D9").Value  = "А7"
$ws.Range("D10").Value = "А8"
$ws.Range("D11").Value = "К1"
$ws.Range("D12").Value = "К2"
$ws.Range("D13").Value = "К3"
$ws.Range("D14").Value = "К4"
$ws.Range("D15").Value = "К5"
$ws.Range("D16").Value = "К6"
$ws.Range("D17").Value = "К7"
$ws.Range("D18").Value = "К8"

# ---- New column E values (control signal from MCU) ----
$ws.Range("E3").Value  = "p"
$ws.Range("E4").Value  = "n"
$ws.Range("E5").Value  = "m"
$ws.Range("E6").Value  = "l"
$ws.Range("E7").Value  = "k"
$ws.Range("E8").Value  = "j"
$ws.Range("E9").Value  = "h"
$ws.Range("E10").Value = "g"
$ws.Range("E11").Value = "a"
$ws.Range("E12").Value = "b"
$ws.Range("E13").Value = "c"
$ws.Range("E14").Value = "d"
$ws.Range("E15").Value = "e"
$ws.Range("E16").Value = "f"
$ws.Range("E17").Value = "dp"
$ws.Range("E18").Value = "nMO"

# ---- Re-use the existing bordered/centered style (copy format only) ----
$ws.Range("B2").Copy()
$ws.Range("D2:E18").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Header row: wrap + taller row for the long "control signal" caption
$ws.Range("E2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 28.8

# ---- Column widths for the new columns ----
$ws.Columns.Item(4).ColumnWidth = 11.33203125
$ws.Columns.Item(5).ColumnWidth = 13.77734375

# ---- Picture: move/resize/rotate ----
$shp = $ws.Shapes.Item(1)
$shp.Rotation = -60
$shp.Left = 2388600 / 914400 * 72
$shp.Top = 161241 / 914400 * 72

# ---- Move selection like in the authored file ----
$ws.Range("N8").Select()
